# Auto-generated Excel COM script applying numeric corrections to the
# FFXIV Leve profitability sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each block below updates currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H-N) for specific rows as produced by the scheduled data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 1997.75
$ws.Range("I51").Value = 2000
$ws.Range("K51").Value = 2000
$ws.Range("M51").Value = -1516

$ws.Range("H62").Value = 25006460
$ws.Range("I62").Value = 1541.6471
$ws.Range("J62").Value = 166701000
$ws.Range("K62").Value = 1541.6471
$ws.Range("L62").Value = 166701000
$ws.Range("M62").Value = -917.6470999999999
$ws.Range("N62").Value = -166702248

$ws.Range("H65").Value = 25006460
$ws.Range("I65").Value = 1541.6471
$ws.Range("J65").Value = 166701000
$ws.Range("K65").Value = 7708.2355
$ws.Range("L65").Value = 833505000
$ws.Range("M65").Value = -4588.2355
$ws.Range("N65").Value = -833511240

$ws.Range("H86").Value = 7424.15
$ws.Range("I86").Value = 8136.2
$ws.Range("K86").Value = 8136.2
$ws.Range("M86").Value = -7013.2

$ws.Range("H89").Value = 7424.15
$ws.Range("I89").Value = 8136.2
$ws.Range("K89").Value = 40681
$ws.Range("M89").Value = -35065

$ws.Range("H116").Value = 6715.294
$ws.Range("I116").Value = 6328.8887
$ws.Range("K116").Value = 6328.8887
$ws.Range("M116").Value = -2886.8887

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H131").Value = 1167.1765
$ws.Range("I131").Value = 988.7857
$ws.Range("J131").Value = 1999.6666
$ws.Range("K131").Value = 2966.3571
$ws.Range("L131").Value = 5998.9998
$ws.Range("M131").Value = 2073.6429
$ws.Range("N131").Value = -16078.9998

$ws.Range("H132").Value = 5210150
$ws.Range("I132").Value = 6025374.5
$ws.Range("J132").Value = 5254.385
$ws.Range("K132").Value = 18076123.5
$ws.Range("L132").Value = 15763.155
$ws.Range("M132").Value = -18073593.5
$ws.Range("N132").Value = -20823.155

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26407.037
$ws.Range("I32").Value = 10756.314
$ws.Range("K32").Value = 10756.314
$ws.Range("M32").Value = -10469.314

$ws.Range("H61").Value = 2283.0557
$ws.Range("I61").Value = 2441
$ws.Range("J61").Value = 1493.3334
$ws.Range("K61").Value = 2441
$ws.Range("L61").Value = 1493.3334
$ws.Range("M61").Value = -2229
$ws.Range("N61").Value = -1917.3334

$ws.Range("H74").Value = 1736.1482
$ws.Range("I74").Value = 1041.7368
$ws.Range("J74").Value = 3385.375
$ws.Range("K74").Value = 1041.7368
$ws.Range("L74").Value = 3385.375
$ws.Range("M74").Value = -167.7367999999999
$ws.Range("N74").Value = -5133.375

$ws.Range("H77").Value = 1736.1482
$ws.Range("I77").Value = 1041.7368
$ws.Range("J77").Value = 3385.375
$ws.Range("K77").Value = 5208.683999999999
$ws.Range("L77").Value = 16926.875
$ws.Range("M77").Value = -840.6839999999993
$ws.Range("N77").Value = -25662.875

$ws.Range("H132").Value = 2914.3948
$ws.Range("I132").Value = 1689.3939
$ws.Range("J132").Value = 10999.4
$ws.Range("K132").Value = 5068.1817
$ws.Range("L132").Value = 32998.2
$ws.Range("M132").Value = -2538.1817
$ws.Range("N132").Value = -38058.2

$ws.Range("H136").Value = 2283.0557
$ws.Range("I136").Value = 2441
$ws.Range("J136").Value = 1493.3334
$ws.Range("K136").Value = 7323
$ws.Range("L136").Value = 4480.0002
$ws.Range("M136").Value = -4773
$ws.Range("N136").Value = -9580.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 12952
$ws.Range("I96").Value = 12952
$ws.Range("K96").Value = 12952
$ws.Range("M96").Value = -10206

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1500
$ws.Range("J2").Value = 3000
$ws.Range("L2").Value = 3000
$ws.Range("N2").Value = -3226

$ws.Range("H31").Value = 3230.7368
$ws.Range("I31").Value = 1738.2273
$ws.Range("J31").Value = 8282.308000000001
$ws.Range("K31").Value = 1738.2273
$ws.Range("L31").Value = 8282.308000000001
$ws.Range("M31").Value = -1443.2273
$ws.Range("N31").Value = -8872.308000000001

$ws.Range("H34").Value = 3230.7368
$ws.Range("I34").Value = 1738.2273
$ws.Range("J34").Value = 8282.308000000001
$ws.Range("K34").Value = 1738.2273
$ws.Range("L34").Value = 8282.308000000001
$ws.Range("M34").Value = -1536.2273
$ws.Range("N34").Value = -8686.308000000001

$ws.Range("H99").Value = 1353
$ws.Range("I99").Value = 1307.4615
$ws.Range("J99").Value = 1412.2
$ws.Range("K99").Value = 1307.4615
$ws.Range("L99").Value = 1412.2
$ws.Range("M99").Value = 190.5385000000001
$ws.Range("N99").Value = -4408.2

$ws.Range("H126").Value = 1353
$ws.Range("I126").Value = 1307.4615
$ws.Range("J126").Value = 1412.2
$ws.Range("K126").Value = 3922.3845
$ws.Range("L126").Value = 4236.6
$ws.Range("M126").Value = -1452.3845
$ws.Range("N126").Value = -9176.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1593.0714
$ws.Range("J11").Value = 1900.2727
$ws.Range("L11").Value = 5700.8181
$ws.Range("N11").Value = -5980.8181

$ws.Range("H42").Value = 3000
$ws.Range("J42").Value = 3000
$ws.Range("L42").Value = 9000
$ws.Range("N42").Value = -10068

$ws.Range("H58").Value = 3292.15
$ws.Range("I58").Value = 1250
$ws.Range("J58").Value = 3519.0557
$ws.Range("K58").Value = 3750
$ws.Range("L58").Value = 10557.1671
$ws.Range("M58").Value = -3622
$ws.Range("N58").Value = -10813.1671

$ws.Range("H122").Value = 1277.081
$ws.Range("I122").Value = 598.1667
$ws.Range("J122").Value = 1602.96
$ws.Range("K122").Value = 5383.5003
$ws.Range("L122").Value = 14426.64
$ws.Range("M122").Value = -2933.5003
$ws.Range("N122").Value = -19326.64

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 28663.373
$ws.Range("I70").Value = 35103.176
$ws.Range("J70").Value = 4335.222
$ws.Range("K70").Value = 35103.176
$ws.Range("L70").Value = 4335.222
$ws.Range("M70").Value = -34833.176
$ws.Range("N70").Value = -4875.222

$ws.Range("H73").Value = 28663.373
$ws.Range("I73").Value = 35103.176
$ws.Range("J73").Value = 4335.222
$ws.Range("K73").Value = 35103.176
$ws.Range("L73").Value = 4335.222
$ws.Range("M73").Value = -34167.176
$ws.Range("N73").Value = -6207.222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1311.125
$ws.Range("I22").Value = 999.75
$ws.Range("J22").Value = 1622.5
$ws.Range("K22").Value = 999.75
$ws.Range("L22").Value = 1622.5
$ws.Range("M22").Value = -704.75
$ws.Range("N22").Value = -2212.5

$ws.Range("H27").Value = 1311.125
$ws.Range("I27").Value = 999.75
$ws.Range("J27").Value = 1622.5
$ws.Range("K27").Value = 999.75
$ws.Range("L27").Value = 1622.5
$ws.Range("M27").Value = -892.75
$ws.Range("N27").Value = -1836.5

$ws.Range("H46").Value = 3225
$ws.Range("I46").Value = 1950
$ws.Range("K46").Value = 1950
$ws.Range("M46").Value = -1762

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H132").Value = 2422.25
$ws.Range("I132").Value = 2070.139
$ws.Range("J132").Value = 4006.75
$ws.Range("K132").Value = 6210.417
$ws.Range("L132").Value = 12020.25
$ws.Range("M132").Value = -3680.417
$ws.Range("N132").Value = -17080.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 15338
$ws.Range("I122").Value = 23496.445
$ws.Range("J122").Value = 4848.5713
$ws.Range("K122").Value = 70489.33499999999
$ws.Range("L122").Value = 14545.7139
$ws.Range("M122").Value = -68039.33499999999
$ws.Range("N122").Value = -19445.7139
